$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Датмоа"
$ws.Range("B3").Value = "Ыдлалвлоа"
$ws.Range("C3").Value = "Ылавооаща"
# D3 becomes a digit-only string; prefix with an apostrophe so Excel stores it
# as text (matching t="inlineStr"/string type in the target) instead of
# silently converting it to a number, then reset the cell style back to
# Normal so no visible formatting change is introduced.
$ws.Range("D3").Value = "'3258748536"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "лыовподыраопып"
